# QAExcelCompiler: Simplify TOTAL stats table - remove Completion%, rename Total to Done
#
# On the TOTAL sheet, the "EN TESTER STATS" table (rows 2-15, originally
# columns A:K) drops its "Completion" (B) and "Actual Issues" (C) columns
# (that info already lives in the Category Breakdown table lower on the
# sheet), and the old "Total" header is renamed to "Done" for consistency
# with the DAILY sheet. Every column from D onward shifts left by two.
#
# We do this per-row (not via a whole-column delete) because rows 18+
# (Category Breakdown / Ranking tables) must stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TOTAL")

# ---- Header row (row 2): shift D2:K2 -> B2:I2 (brings the right styles
# along, e.g. F2/G2 pick up the "s=5" look that used to belong to H2/I2),
# then stamp the renamed labels on top. ----
$ws.Range("D2:K2").Copy($ws.Range("B2"))
$ws.Range("B2").Value = "Done"
$ws.Range("C2").Value = "Issues"
$ws.Range("D2").Value = "No Issue"
$ws.Range("E2").Value = "Blocked"
$ws.Range("F2").Value = "Fixed"
$ws.Range("G2").Value = "Reported"
$ws.Range("H2").Value = "Checking"
$ws.Range("I2").Value = "Pending"
# J2/K2 (old "Checking"/"Pending" headers) are now gone
$ws.Range("J2:K2").Clear()

# ---- Data rows: shift D:K left into B:I (values + styles), drop J:K ----
$dataRows = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 12, 15)
foreach ($r in $dataRows) {
    $ws.Range("D$($r):K$($r)").Copy($ws.Range("B$($r)"))
    $ws.Range("J$($r):K$($r)").Clear()
}

# ---- Shrink the title merges from A:K to A:I to match the new width ----
$ws.Range("A1:K1").UnMerge()
$ws.Range("A1:I1").Merge()
$ws.Range("A14:K14").UnMerge()
$ws.Range("A14:I14").Merge()
